$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change I7 formula to multiply the sum by 7
$ws.Range("I7").Formula = "=(+E7+F7+G7+H7)*7"

# Fill the same formula pattern down into I8:I9 (creates a shared formula group)
$ws.Range("I8:I9").Formula = "=(+E8+F8+G8+H8)*7"

# Update the selected cell to J7
$ws.Range("J7").Select()
